# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update Panama statistics (row 49) ---
$ws.Range("B49").Value = 7197
$ws.Range("C49").Value = 107
$ws.Range("D49").Value = 641
$ws.Range("E49").Value = 6356
$ws.Range("F49").Value = 91
$ws.Range("G49").Value = 3
$ws.Range("H49").Value = 200

# --- Reorder country list entries (Santa Lucia / Belice swap, rows 188-189) ---
$ws.Range("A188").Value = "Santa Lucia"
$ws.Range("B188").Value = 18
$ws.Range("C188").Value = 0
$ws.Range("D188").Value = 15
$ws.Range("E188").Value = 3
$ws.Range("F188").Value = 0
$ws.Range("G188").Value = 0
$ws.Range("H188").Value = 0

$ws.Range("A189").Value = "Belice"
$ws.Range("B189").Value = 18
$ws.Range("C189").Value = 0
$ws.Range("D189").Value = 13
$ws.Range("E189").Value = 3
$ws.Range("F189").Value = 1
$ws.Range("G189").Value = 0
$ws.Range("H189").Value = 2

# --- Reorder country list entries (San Vicente / Namibia swap, rows 194-195) ---
$ws.Range("A194").Value = "Namibia"
$ws.Range("A195").Value = "San Vicente y las Granadinas"

# --- Reorder country list entries (San Cristobal y Nieves / Burundi swap, rows 198-199) ---
$ws.Range("A198").Value = "Burundi"
$ws.Range("B198").Value = 15
$ws.Range("C198").Value = 0
$ws.Range("D198").Value = 7
$ws.Range("E198").Value = 7
$ws.Range("F198").Value = 0
$ws.Range("G198").Value = 0
$ws.Range("H198").Value = 1

$ws.Range("A199").Value = "San Cristobal y Nieves"
$ws.Range("B199").Value = 15
$ws.Range("C199").Value = 0
$ws.Range("D199").Value = 8
$ws.Range("E199").Value = 7
$ws.Range("F199").Value = 0
$ws.Range("G199").Value = 0
$ws.Range("H199").Value = 0
